$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column D ("CodigoVendedor") before the existing "FechaModificacion" column,
# shifting the old FechaModificacion data into column E.
$ws.Columns.Item(4).Insert(-4161)

# New header for column D
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "CodigoVendedor"

# New value for row 2 in column D
$ws.Range("D2").Value = "B024"

# Update the password (B2) and the FechaModificacion value (now in E2)
$ws.Range("B2").Value = "gAAAAABnt587EQgIJ1MPMEoc-hBYkmIEyHuyNhxhZKIlMkggewm5uRtSuHG9gpcixYu7gRNw5iEmCzwLgDTBV48lrxC7bWVEAg=="
$ws.Range("E2").Value = "2025-02-20 16:31:39"
